# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" (E18:E43) list for JULIETH RAMIREZ RAMOS is reversed
# (was ascending 2207..2408, now descending 2408..2207), and the two
# "Valor Mora" (F) amounts that were tied to the first/last period
# (67625 and 47337) travel along with their periods (i.e. they are
# swapped between row 18 and row 43).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Current (ascending) order of periods found in E18:E43.
$periods = @("2207","2208","2209","2210","2211","2212", `
             "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312", `
             "2401","2402","2403","2404","2405","2406","2407","2408")

$firstRow = 18
$lastRow = 43

# Reverse the periods so the newest period is on top.
$reversed = @()
for ($i = $periods.Length - 1; $i -ge 0; $i--) {
    $reversed += $periods[$i]
}

for ($i = 0; $i -lt $reversed.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $reversed[$i]
}

# The "Valor Mora" values swap along with the periods: row 18 (now 2408)
# takes the amount that used to sit on row 43, and row 43 (now 2207)
# takes the amount that used to sit on row 18.
$ws.Cells.Item($firstRow, 6).Value = 47337
$ws.Cells.Item($lastRow, 6).Value = 67625
